$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 25.02.2022 11:45"

# Update row 10 (EuroOil Opustena) values
$ws.Range("B10").Value = 37.9
$ws.Range("C10").Value = 37.7

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "+0.2"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2022-02-25 11:47:32"
$ws.Range("E10").Style = "Normal"
